# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) on the leve-profit rows that changed this run, per-job (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 18
$ws.Range("H18").Value = 1237.375
$ws.Range("I18").Value = 799.6667
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 799.6667
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -515.6667
$ws.Range("N18").Value = -2068

# Row 64
$ws.Range("H64").Value = 3267.5833
$ws.Range("I64").Value = 3200
$ws.Range("J64").Value = 3277.238
$ws.Range("K64").Value = 3200
$ws.Range("L64").Value = 3277.238
$ws.Range("M64").Value = -2952
$ws.Range("N64").Value = -3773.238

# Row 67
$ws.Range("H67").Value = 3267.5833
$ws.Range("I67").Value = 3200
$ws.Range("J67").Value = 3277.238
$ws.Range("K67").Value = 3200
$ws.Range("L67").Value = 3277.238
$ws.Range("M67").Value = -2342
$ws.Range("N67").Value = -4993.237999999999

# Row 74
$ws.Range("H74").Value = 3550.5
$ws.Range("I74").Value = 3466.6667
$ws.Range("J74").Value = 3600.8
$ws.Range("K74").Value = 3466.6667
$ws.Range("L74").Value = 3600.8
$ws.Range("M74").Value = -2530.6667
$ws.Range("N74").Value = -5472.8

# Row 77
$ws.Range("H77").Value = 3550.5
$ws.Range("I77").Value = 3466.6667
$ws.Range("J77").Value = 3600.8
$ws.Range("K77").Value = 17333.3335
$ws.Range("L77").Value = 18004
$ws.Range("M77").Value = -12653.3335
$ws.Range("N77").Value = -27364

# Row 100
$ws.Range("H100").Value = 3285.7144
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -5082

# Row 112
$ws.Range("H112").Value = 1944.2858
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1944.2858
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 5832.857400000001
$ws.Range("N112").Value = -8048.857400000001
$ws.Range("M112").ClearContents()


$ws = $wb.Worksheets.Item("ARM")

# Row 63
$ws.Range("H63").Value = 2578.6843
$ws.Range("I63").Value = 2221.6667
$ws.Range("J63").Value = 2900
$ws.Range("K63").Value = 2221.6667
$ws.Range("L63").Value = 2900
$ws.Range("M63").Value = -1535.6667
$ws.Range("N63").Value = -4272

# Row 66
$ws.Range("H66").Value = 2578.6843
$ws.Range("I66").Value = 2221.6667
$ws.Range("J66").Value = 2900
$ws.Range("K66").Value = 11108.3335
$ws.Range("L66").Value = 14500
$ws.Range("M66").Value = -7676.333500000001
$ws.Range("N66").Value = -21364

# Row 80
$ws.Range("H80").Value = 22450
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 22450
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 22450
$ws.Range("N80").Value = -24446
$ws.Range("M80").ClearContents()

# Row 82
$ws.Range("H82").Value = 20666.666
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 29500
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 29500
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -30222

# Row 83
$ws.Range("H83").Value = 22450
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 22450
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 67350
$ws.Range("N83").Value = -77334
$ws.Range("M83").ClearContents()

# Row 85
$ws.Range("H85").Value = 20666.666
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 29500
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 29500
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -31996


$ws = $wb.Worksheets.Item("BSM")

# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# Row 86
$ws.Range("H86").Value = 2024.2
$ws.Range("I86").Value = 1947.2
$ws.Range("J86").Value = 2139.7
$ws.Range("K86").Value = 1947.2
$ws.Range("L86").Value = 2139.7
$ws.Range("M86").Value = -824.2
$ws.Range("N86").Value = -4385.7

# Row 89
$ws.Range("H89").Value = 2024.2
$ws.Range("I89").Value = 1947.2
$ws.Range("J89").Value = 2139.7
$ws.Range("K89").Value = 9736
$ws.Range("L89").Value = 10698.5
$ws.Range("M89").Value = -4120
$ws.Range("N89").Value = -21930.5

# Row 110
$ws.Range("H110").Value = 26700
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 26700
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 26700
$ws.Range("N110").Value = -34880


$ws = $wb.Worksheets.Item("CRP")

# Row 50
$ws.Range("H50").Value = 12999.667
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 12999.667
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 12999.667
$ws.Range("N50").Value = -14249.667

# Row 124
$ws.Range("H124").Value = 11140
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 11140
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 11140
$ws.Range("N124").Value = -16050

# Row 141
$ws.Range("H141").Value = 33789.715
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 33789.715
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 33789.715
$ws.Range("N141").Value = -44149.715


$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 4919.1665
$ws.Range("I3").Value = 2505
$ws.Range("J3").Value = 7333.3335
$ws.Range("K3").Value = 7515
$ws.Range("L3").Value = 22000.0005
$ws.Range("M3").Value = -7403
$ws.Range("N3").Value = -22224.0005

# Row 13
$ws.Range("H13").Value = 40
$ws.Range("I13").Value = 40
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 120
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 48

# Row 64
$ws.Range("H64").Value = 6365.875
$ws.Range("I64").Value = 4999.5
$ws.Range("J64").Value = 6821.3335
$ws.Range("K64").Value = 14998.5
$ws.Range("L64").Value = 20464.0005
$ws.Range("M64").Value = -14728.5
$ws.Range("N64").Value = -21004.0005

# Row 67
$ws.Range("H67").Value = 6365.875
$ws.Range("I67").Value = 4999.5
$ws.Range("J67").Value = 6821.3335
$ws.Range("K67").Value = 14998.5
$ws.Range("L67").Value = 20464.0005
$ws.Range("M67").Value = -14062.5
$ws.Range("N67").Value = -22336.0005

# Row 131
$ws.Range("H131").Value = 1188.9512
$ws.Range("I131").Value = 1998.2858
$ws.Range("J131").Value = 1113.4133
$ws.Range("K131").Value = 5994.857400000001
$ws.Range("L131").Value = 3340.2399
$ws.Range("M131").Value = -954.8574000000008
$ws.Range("N131").Value = -13420.2399


$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 5286.0454
$ws.Range("I70").Value = 5163.706
$ws.Range("J70").Value = 5702
$ws.Range("K70").Value = 5163.706
$ws.Range("L70").Value = 5702
$ws.Range("M70").Value = -4893.706
$ws.Range("N70").Value = -6242

# Row 73
$ws.Range("H73").Value = 5286.0454
$ws.Range("I73").Value = 5163.706
$ws.Range("J73").Value = 5702
$ws.Range("K73").Value = 5163.706
$ws.Range("L73").Value = 5702
$ws.Range("M73").Value = -4227.706
$ws.Range("N73").Value = -7574

# Row 110
$ws.Range("H110").Value = 32875
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 32875
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 32875
$ws.Range("N110").Value = -41055


$ws = $wb.Worksheets.Item("LTW")

# Row 41
$ws.Range("H41").Value = 26000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 26000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 26000
$ws.Range("N41").Value = -26876

# Row 110
$ws.Range("H110").Value = 26469.715
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 26469.715
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 26469.715
$ws.Range("N110").Value = -34649.715

# Row 136
$ws.Range("H136").Value = 5400.76
$ws.Range("I136").Value = 1247.5883
$ws.Range("J136").Value = 14226.25
$ws.Range("K136").Value = 3742.7649
$ws.Range("L136").Value = 42678.75
$ws.Range("M136").Value = -1192.7649
$ws.Range("N136").Value = -47778.75


$ws = $wb.Worksheets.Item("WVR")

# Row 58
$ws.Range("H58").Value = 11398.8
$ws.Range("I58").Value = 9000
$ws.Range("J58").Value = 12998
$ws.Range("K58").Value = 9000
$ws.Range("L58").Value = 12998
$ws.Range("M58").Value = -8692
$ws.Range("N58").Value = -13614

# Row 86
$ws.Range("H86").Value = 27580
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 27580
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 27580
$ws.Range("N86").Value = -29826

# Row 89
$ws.Range("H89").Value = 27580
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 27580
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 137900
$ws.Range("N89").Value = -149132

# Row 99
$ws.Range("H99").Value = 28738.5
$ws.Range("I99").Value = 29215.5
$ws.Range("J99").Value = 28500
$ws.Range("K99").Value = 29215.5
$ws.Range("L99").Value = 28500
$ws.Range("M99").Value = -26220.5
$ws.Range("N99").Value = -34490

